$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet from "1" to "შუახევი"
$ws.Name = "შუახევი"

# Delete old row 2 "(მოსახლეობის აღწერის შედეგებით)" note, and merge the two-year columns
# down into a single 2014 column, shifting rows up:
#   old: Row1 title | Row2 note | Row3 blank | Row4 (კვ. კმ) | Row5 years(1989,2002,2014) | Row6 ფართობი(588,588,588)
#   new: Row1 title | Row2 blank | Row3 (კვ. კმ) | Row4 blank,2014 | Row5 ფართობი,588

# Remove the note cell in A2 (clear content + formatting)
$ws.Range("A2").Clear()

# Delete old row 3 (was a spacer row with only a formatted empty B3) - shifts rows 4,5,6 up to 3,4,5
$ws.Rows("3").Delete()

# Delete columns B:C so that only the 2014 value (was column D) remains, shifting it to column B
$ws.Range("B:C").Delete()

# The "ფართობი" label cell (now A5) loses its left-edge medium border now that the
# table is a single data column wide
$ws.Range("A5").Borders(7).LineStyle = 0   # xlEdgeLeft = 7, xlLineStyleNone = 0

# Leave the selection on A2, matching the saved view state
$ws.Range("A2").Select() | Out-Null

$wb.Save()
